# Update cryptocurrency price / volume figures per the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "65.455.51"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "3.560.92"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextCell "D5" "600.49"
$ws.Range("E5").Value = "  +0.53%  "
Set-TextCell "D6" "140.54"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("D7").Value = "3.560.74"
$ws.Range("E7").Value = "  +2.98%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +2.65%  "
$ws.Range("E11").Value = "  -6.12%  "
Set-TextCell "D12" "0.396"
$ws.Range("E12").Value = "  +3.82%  "
$ws.Range("D13").Value = "4.163.44"
$ws.Range("E13").Value = "  +3.08%  "
$ws.Range("E14").Value = "  +2.38%  "
Set-TextCell "D15" "27.16"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "3.562.31"
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "65.517.83"
$ws.Range("E18").Value = "  -0.27%  "
Set-TextCell "D19" "10.33"
$ws.Range("E19").Value = "  +4.03%  "
Set-TextCell "D20" "5.90"
$ws.Range("E20").Value = "  +1.48%  "
Set-TextCell "D21" "14.26"
$ws.Range("E21").Value = "  +3.34%  "
Set-TextCell "D22" "396.67"
$ws.Range("E22").Value = "  +0.16%  "
Set-TextCell "D23" "0.574"
$ws.Range("E23").Value = "  +4.41%  "
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.703.66"
$ws.Range("E24").Value = "  +2.87%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D25" "74.39"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +8.80%  "
Set-TextCell "D28" "7.90"
$ws.Range("E28").Value = "  +8.29%  "
Set-TextCell "D29" "0.999"
$ws.Range("E29").Value = "  -0.02%  "
Set-TextCell "D30" "2.29"
$ws.Range("E30").Value = "  +0.18%  "
Set-TextCell "D31" "8.31"
$ws.Range("D32").Value = "3.576.82"
$ws.Range("E32").Value = "  +3.37%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell "D33" "1.00"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D34" "23.93"
$ws.Range("E34").Value = "  +4.06%  "
Set-TextCell "D35" "0.147"
$ws.Range("E35").Value = "  +0.90%  "
Set-TextCell "D36" "1.27"
$ws.Range("E36").Value = "  +3.29%  "
Set-TextCell "D37" "7.08"
$ws.Range("E37").Value = "  +1.50%  "
Set-TextCell "D38" "168.62"
$ws.Range("E38").Value = "  -2.71%  "
Set-TextCell "D39" "1.55"
$ws.Range("E39").Value = "  +0.71%  "
Set-TextCell "D40" "5.03"
$ws.Range("E40").Value = "  +3.67%  "
Set-TextCell "D41" "0.0807"
$ws.Range("E41").Value = "  +3.22%  "
Set-TextCell "D42" "0.835"
$ws.Range("E42").Value = "  +1.50%  "
Set-TextCell "D43" "26.88"
$ws.Range("E43").Value = "  +15.11%  "
Set-TextCell "D44" "42.89"
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("E45").Value = "  -0.04%  "
Set-TextCell "D46" "4.45"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("E47").Value = "  +3.06%  "
Set-TextCell "D48" "1.20"
$ws.Range("E48").Value = "  +7.07%  "
$ws.Range("D49").Value = "2.442.57"
$ws.Range("E49").Value = "  +10.33%  "
Set-TextCell "D50" "6.83"
$ws.Range("E50").Value = "  +3.53%  "
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D51" "2.38"
$ws.Range("E51").Value = "  +19.98%  "
